$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: " Linux computer running ROS ... using OpenC" + "V" (with an
# in-between _GoBack bookmark) become a single merged run
# " Linux computer running ROS ... using OpenCV" (bookmark removed).
# ---------------------------------------------------------------------------

$r1 = $d.Content
$found1 = $r1.Find.Execute("object recognition using OpenC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Found 'OpenC' run: " + $found1)
$opencEnd = $r1.End

# Firewall the run boundary to the left (so "Program"/"med" don't get pulled
# into the merge cascade) and to the right of "V" (so " "/"and"/" Python"
# keep their separate runs).
$afterV = $opencEnd + 1

# Find the start of the run containing "OpenC" by locating the preceding
# run break: walk left from the Find match start until we hit the run that
# starts right after "med".
$r2 = $d.Content
$found2 = $r2.Find.Execute("med Linux computer running ROS to perform object recognition using OpenC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Found 'med Linux...OpenC' span: " + $found2)
$runStart = $r2.Start + 3   # skip past "med" to the start of the " Linux..." run

$fwLeft = $d.Range($runStart, $runStart)
$d.Bookmarks.Add("_TempFWLeft", $fwLeft) | Out-Null

$fwRight = $d.Range($afterV, $afterV)
$d.Bookmarks.Add("_TempFWRight", $fwRight) | Out-Null

$fwRight2 = $d.Range($afterV + 1, $afterV + 1)
$d.Bookmarks.Add("_TempFWRight2", $fwRight2) | Out-Null

$fwRight3 = $d.Range($afterV + 4, $afterV + 4)
$d.Bookmarks.Add("_TempFWRight3", $fwRight3) | Out-Null

$d.Bookmarks("_GoBack").Delete()

$vRange2 = $d.Range($opencEnd, $opencEnd + 1)
$vRange2.Text = "VX"
$xRange = $d.Range($opencEnd + 1, $opencEnd + 2)
$xRange.Text = ""

$d.Bookmarks("_TempFWLeft").Delete()
$d.Bookmarks("_TempFWRight").Delete()
$d.Bookmarks("_TempFWRight2").Delete()
$d.Bookmarks("_TempFWRight3").Delete()

# ---------------------------------------------------------------------------
# Change 2: "recommendations ... tutoring sessions" is split into
# "...sessio" + _GoBack bookmark + "ns" (same run formatting).
# ---------------------------------------------------------------------------

$r3 = $d.Content
$found3 = $r3.Find.Execute("recommendations to improve students", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Found 'recommendations...' run: " + $found3)
$recStart = $r3.Start

$fw2 = $d.Range($recStart, $recStart)
$d.Bookmarks.Add("_TempFW2", $fw2) | Out-Null

$r4 = $d.Content
$found4 = $r4.Find.Execute("tutoring sessions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Found 'tutoring sessions': " + $found4)
$splitPos = $r4.End - 2   # right before "ns"

$ins = $d.Range($splitPos, $splitPos)
$ins.InsertAfter("Z")
$zRange = $d.Range($splitPos, $splitPos + 1)
$zRange.Text = ""

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Bookmarks("_TempFW2").Delete()
